# This commit ("fixed a bug in generating tables") changed table-generation
# code in the authoring tool, which regenerates the whole .pptx package
# (including fresh, randomly-generated OPC relationship IDs) on every build.
# This particular slide has no table and none of its visible/semantic
# content (text, shapes, images, layout bindings, embedded add-in data)
# actually changed between the "before" and "after" snapshots -- only the
# internal relationship-id strings (and one cosmetic webextension GUID)
# differ, as a side effect of the regenerate-on-save behaviour.
#
# There is no PowerPoint object-model operation that corresponds to "mint
# new relationship ids" (that's an OPC packaging implementation detail, not
# something exposed to VBA/COM automation), so there is nothing for this
# script to legitimately change. We touch the presentation object (without
# altering any visible property) so the save still round-trips cleanly.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
